$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("E3").Value = "2016-03-17 18:34:35"
$wsZh.Range("H3").Value = "2016-03-17 18:34:54"

$wsDe.Range("E3").Value = "2016-03-17 18:34:39"
$wsDe.Range("H3").Value = "2016-03-17 18:35:00"
